# Applies the PlayerPerformance_3978 update:
#  1. Clears the (already-empty) INNING_NUMBER cells B4/B6/B8/B10 on
#     "ODI Batting" so the <c> element disappears entirely, matching rows
#     where the player did not bat.
#  2. Adds a new "ODI Batting Extra" sheet (after "ODI Bowling") with
#     per-match batting-position / boundary-count / man-of-the-match data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "ODI Batting" - drop the stray empty INNING_NUMBER cells.
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B4").ClearContents()
$batting.Range("B6").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B10").ClearContents()

# ---------------------------------------------------------------------
# 2. New sheet "ODI Batting Extra", appended after the last sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Reuse the bold/bordered header style already used by the other sheets
# instead of re-building the font/border/alignment combo from scratch.
$headerSrc = $batting.Range("A1:F1")
$headerSrc.Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $extra.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# Helper: write a value as literal text even when it looks like a number
# or a percentage (Excel would otherwise silently coerce "4100"/"6.02%"
# into a numeric/percent cell).
function Set-TextValue {
    param($cell, $value)
    if ($value -eq $null) {
        return
    }
    if ($value -match '^[0-9]+(\.[0-9]+)?%?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4100", 7,    "2",  "0",  "6.02%",  "NO"),
    @("4101", 7,    "1",  "2",  "15.08%", "NO"),
    @("4102", $null, $null, $null, $null, "NO"),
    @("4103", 9,    "0",  "0",  $null,    "NO"),
    @("4104", 8,    $null, $null, $null,  "NO"),
    @("4105", 8,    "1",  "0",  "1.95%",  "NO"),
    @("4248", 9,    $null, $null, $null,  "NO"),
    @("4249", 9,    "1",  "0",  "4.61%",  "NO"),
    @("4251", 9,    $null, $null, $null,  "NO")
)

$r = 2
foreach ($row in $rows) {
    Set-TextValue $extra.Cells.Item($r, 1) $row[0]
    $battingPos = $row[1]
    if ($battingPos -ne $null) {
        $extra.Cells.Item($r, 2).Value = $battingPos
    }
    Set-TextValue $extra.Cells.Item($r, 3) $row[2]
    Set-TextValue $extra.Cells.Item($r, 4) $row[3]
    Set-TextValue $extra.Cells.Item($r, 5) $row[4]
    Set-TextValue $extra.Cells.Item($r, 6) $row[5]
    $r = $r + 1
}

# Restore the originally-active sheet/tab now that the new sheet has been
# built (adding/naming it shifts Excel's active-sheet focus to it).
$wb.Worksheets.Item(1).Activate()

